$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ligue 1"
$ws.Range("B2").Value = "Paradou AC"
$ws.Range("C2").Value = "Belouizdad"
$ws.Range("D2").Value = 2.65
$ws.Range("E2").Value = 2.4
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 36
